$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.319.16"
$ws.Range("E2").Value = "  -0.31%  "

$ws.Range("D3").Value = "1.841.11"
$ws.Range("E3").Value = "  -0.67%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9993"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "239.98"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.42%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6265"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.20%  "

$ws.Range("E7").Value = "  +0.06%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07406"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.34%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2891"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.35%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.72"
$ws.Range("D10").Style = "Normal"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07737"
$ws.Range("D11").Style = "Normal"

$ws.Range("D12").Value = "1.838.56"
$ws.Range("E12").Value = "  -0.79%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.975"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.02%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6762"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.79%  "

$ws.Range("E15").Value = "  -2.94%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "81.96"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.69%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.230"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.40%  "

$ws.Range("D18").Value = "29.296.34"
$ws.Range("E18").Value = "  -0.42%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "228.55"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.81%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.28"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.92%  "

$ws.Range("E21").Value = "  +0.03%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.409"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.16%  "

$ws.Range("E23").Value = "  +0.17%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "158.68"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.30%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "8.457"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.05%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1347"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.65%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.38"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.63%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.06655"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +16.87%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.449"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.51%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.483"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.54%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.059"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.77%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.061"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.13%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.831"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.06%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.136"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.93%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6917"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.31%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.571"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.50%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01856"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.52%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.827"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.59%  "

$ws.Range("D39").Value = "1.243.07"
$ws.Range("E39").Value = "  -0.35%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.745"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.37%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9345"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.40%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.000"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.02%  "

$ws.Range("D43").Value = "1.981.86"
$ws.Range("E43").Value = "  -1.48%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "100.66"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.36%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "65.31"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.10%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.024"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.77%  "

$ws.Range("E47").Value = "  +1.49%  "

$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.00000000116"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.53%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.000"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.20%  "

$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.1149"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.70%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3888"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.01%  "
